# Edit script: update four marketing slogan/tagline sentences in the
# Portuguese document to match the target copy (Juno / OLPRODLOC check-in).
#
# NOTE: Find.Execute's built-in Replace mode smart-quotes straight
# quotes/apostrophes in the replacement text, which would corrupt the
# "straight quote" style already used throughout this document. To avoid
# that, we locate each target sentence with Find.Execute (no replace) and
# then assign the new sentence directly to the found Range's .Text, which
# performs a literal (non-autoformatted) substitution.

$d = $word.ActiveDocument

function Replace-Exact($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute(
        $find,      # FindText
        $true,      # MatchCase
        $true,      # MatchWholeWord
        $false,     # MatchWildcards
        $false,     # MatchSoundsLike
        $false,     # MatchAllWordForms
        $true,      # Forward
        1,          # Wrap (wdFindContinue)
        $false,     # Format
        "",         # ReplaceWith (unused - we set Range.Text ourselves)
        0           # Replace (wdReplaceNone)
    )
    if (-not $found) {
        throw "Find failed for: $find"
    }
    $range.Text = $replace
}

Replace-Exact `
    "A campanha de marketing usará o seguinte slogan para capturar a essência da marca Munson: `"Munson's: Pickles and Ppreserve with a Purpose`"." `
    "A campanha de marketing usará o seguinte slogan para capturar a essência da marca Munson's: `"Munson's: Pickles and Preserves com propósito`"."

Replace-Exact `
    "A campanha de marketing usará o seguinte slogan para enfatizar os benefícios do produto Munson: `"Munson's: More than Just Pickles and Preserves`"." `
    "A campanha de marketing usará o seguinte slogan para enfatizar os benefícios do produto da Munson: `"Munson's: Mais do que apenas Pickles and Preserves`"."

Replace-Exact `
    "A campanha de marketing usará o seguinte lema para inspirar a defesa do cliente de Munson: `"Munson's: Share the Love of Pickles and Preserves`"." `
    "A campanha de marketing usará o seguinte lema para inspirar a defesa do cliente da Munson's: `"Munson's: Compartilhe o amor por Pickles and Preserves`"."

Replace-Exact `
    "A campanha de marketing usará a seguinte frase para impulsionar o teste e a compra do produto de Munson: `"Munson's: Find Them, Try Them, Love Them`"." `
    "A campanha de marketing usará a seguinte frase para impulsionar a experimentação e a compra do produto da Munson's: `"Munson's: encontre, experimente, apaixone-se`"."
